$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the historical date values (fracciones / historico)
$ws.Range("B8").Value = 44743
$ws.Range("C8").Value = 44834
$ws.Range("F8").Value = 44844
$ws.Range("G8").Value = 44844

# Update the active cell / selection saved in the worksheet view
$ws.Range("B11").Select()
